$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected Hydrogen figures (row 3) ---
# Overall Iron & steel demand revised upward
$ws.Range("B3").Value = 17885007.67964587
# Non-metallic minerals Hydrogen figure removed (no longer applicable)
# (set via a quote-prefixed empty string, then reset the style, so the
#  cell keeps its blank-text identity like its neighbouring blank cells
#  instead of being dropped from the sheet entirely)
$ws.Range("D3").Value = "'"
$ws.Range("D3").Style = "Normal"

# --- Corrected Methanol figure (row 4, Chemicals) ---
$ws.Range("C4").Value = 43.40194001924809

# --- Corrected Ammonia figure (row 5, Chemicals) ---
$ws.Range("C5").Value = 3101.662712617846

# --- Row 7 relabelled from "Other" to "Biogas" with corrected figure ---
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 4539.65494025398

# --- New row 8: "Other" (split out of the old "Other" row) ---
# Match the row-label formatting used by the other category cells in column A
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A8").Value = "Other"

$ws.Range("B8").Value = "'"
$ws.Range("B8").Style = "Normal"

$ws.Range("C8").Value = "'"
$ws.Range("C8").Style = "Normal"

$ws.Range("D8").Value = 848.7396134232458
